$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8116
$ws1.Range("F7").Value = 2053
$ws1.Range("F9").Value = 43
$ws1.Range("F15").Value = 8474
$ws1.Range("F30").Value = 2058
$ws1.Range("F36").Value = 142

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 388

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2320

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2320
$ws4.Range("F4").Value = 388
$ws4.Range("F6").Value = 8116
$ws4.Range("F11").Value = 2053
$ws4.Range("F14").Value = 43
$ws4.Range("F23").Value = 8474
$ws4.Range("F35").Value = 2058
$ws4.Range("F41").Value = 142
